$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 911.44446
$ws.Range("I2").Value = 525.375
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 525.375
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -412.375
$ws.Range("N2").Value = -4226

$ws.Range("H43").Value = 1324.75
$ws.Range("I43").Value = 700
$ws.Range("J43").Value = 1414
$ws.Range("K43").Value = 700
$ws.Range("L43").Value = 1414
$ws.Range("M43").Value = -631
$ws.Range("N43").Value = -1552

$ws.Range("H69").Value = 5257.5
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 5294.2856
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 15882.8568
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -17630.8568

$ws.Range("H72").Value = 5257.5
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 5294.2856
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 47648.5704
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -56384.5704

$ws.Range("H132").Value = 2637.7693
$ws.Range("I132").Value = 2699.25
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 8097.75
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -5567.75
$ws.Range("N132").Value = -10760

$ws.Range("H133").Value = 34508
$ws.Range("J133").Value = 34508
$ws.Range("L133").Value = 34508
$ws.Range("N133").Value = -44628

$ws.Range("H137").Value = 1378.2167
$ws.Range("I137").Value = 1336.341
$ws.Range("J137").Value = 1493.375
$ws.Range("K137").Value = 4009.023
$ws.Range("L137").Value = 4480.125
$ws.Range("M137").Value = -1459.023
$ws.Range("N137").Value = -9580.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7260.9297
$ws.Range("I32").Value = 4511.3438
$ws.Range("K32").Value = 4511.3438
$ws.Range("M32").Value = -4224.3438

$ws.Range("H61").Value = 6083.778
$ws.Range("I61").Value = 6919.1
$ws.Range("J61").Value = 3697.1428
$ws.Range("K61").Value = 6919.1
$ws.Range("L61").Value = 3697.1428
$ws.Range("M61").Value = -6707.1
$ws.Range("N61").Value = -4121.1428

$ws.Range("H122").Value = 1116400.2
$ws.Range("I122").Value = 1351167.4
$ws.Range("J122").Value = 1257
$ws.Range("K122").Value = 4053502.2
$ws.Range("L122").Value = 3771
$ws.Range("M122").Value = -4051052.2
$ws.Range("N122").Value = -8671

$ws.Range("H132").Value = 2640.0588
$ws.Range("I132").Value = 1750.6086
$ws.Range("J132").Value = 4499.8184
$ws.Range("K132").Value = 5251.825800000001
$ws.Range("L132").Value = 13499.4552
$ws.Range("M132").Value = -2721.825800000001
$ws.Range("N132").Value = -18559.4552

$ws.Range("H136").Value = 6083.778
$ws.Range("I136").Value = 6919.1
$ws.Range("J136").Value = 3697.1428
$ws.Range("K136").Value = 20757.3
$ws.Range("L136").Value = 11091.4284
$ws.Range("M136").Value = -18207.3
$ws.Range("N136").Value = -16191.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1096.4783
$ws.Range("I107").Value = 951.3077
$ws.Range("K107").Value = 951.3077
$ws.Range("M107").Value = 968.6923

$ws.Range("H134").Value = 8199.368
$ws.Range("I134").Value = 11982.818
$ws.Range("J134").Value = 2997.125
$ws.Range("K134").Value = 35948.454
$ws.Range("L134").Value = 8991.375
$ws.Range("M134").Value = -33413.454
$ws.Range("N134").Value = -14061.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4995.085
$ws.Range("I31").Value = 1792.2572
$ws.Range("J31").Value = 14336.667
$ws.Range("K31").Value = 1792.2572
$ws.Range("L31").Value = 14336.667
$ws.Range("M31").Value = -1497.2572
$ws.Range("N31").Value = -14926.667

$ws.Range("H34").Value = 4995.085
$ws.Range("I34").Value = 1792.2572
$ws.Range("J34").Value = 14336.667
$ws.Range("K34").Value = 1792.2572
$ws.Range("L34").Value = 14336.667
$ws.Range("M34").Value = -1590.2572
$ws.Range("N34").Value = -14740.667

$ws.Range("H132").Value = 2245.0952
$ws.Range("I132").Value = 1679.5834
$ws.Range("J132").Value = 2999.111
$ws.Range("K132").Value = 5038.7502
$ws.Range("L132").Value = 8997.332999999999
$ws.Range("M132").Value = -2508.7502
$ws.Range("N132").Value = -14057.333

$ws.Range("H134").Value = 3246.75
$ws.Range("I134").Value = 3303.7585
$ws.Range("J134").Value = 3010.5715
$ws.Range("K134").Value = 9911.2755
$ws.Range("L134").Value = 9031.7145
$ws.Range("M134").Value = -7376.2755
$ws.Range("N134").Value = -14101.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 307323.88
$ws.Range("J4").Value = 2411.5386
$ws.Range("L4").Value = 7234.6158
$ws.Range("N4").Value = -7458.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3242.4666
$ws.Range("I132").Value = 3675
$ws.Range("J132").Value = 3085.182
$ws.Range("K132").Value = 11025
$ws.Range("L132").Value = 9255.545999999998
$ws.Range("M132").Value = -8495
$ws.Range("N132").Value = -14315.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 101930.3
$ws.Range("I7").Value = 101930.3
$ws.Range("K7").Value = 101930.3
$ws.Range("M7").Value = -101818.3

$ws.Range("H82").Value = 461311.53
$ws.Range("I82").Value = 1251787.1
$ws.Range("J82").Value = 66073.75
$ws.Range("K82").Value = 1251787.1
$ws.Range("L82").Value = 66073.75
$ws.Range("M82").Value = -1251426.1
$ws.Range("N82").Value = -66795.75

$ws.Range("H85").Value = 461311.53
$ws.Range("I85").Value = 1251787.1
$ws.Range("J85").Value = 66073.75
$ws.Range("K85").Value = 1251787.1
$ws.Range("L85").Value = 66073.75
$ws.Range("M85").Value = -1250539.1
$ws.Range("N85").Value = -68569.75

$ws.Range("H122").Value = 2469914.8
$ws.Range("I122").Value = 3107707.5
$ws.Range("J122").Value = 1002991.5
$ws.Range("K122").Value = 9323122.5
$ws.Range("L122").Value = 3008974.5
$ws.Range("M122").Value = -9320672.5
$ws.Range("N122").Value = -3013874.5

$ws.Range("H126").Value = 101930.3
$ws.Range("I126").Value = 101930.3
$ws.Range("K126").Value = 305790.9
$ws.Range("M126").Value = -303320.9

$ws.Range("H132").Value = 6969496.5
$ws.Range("I132").Value = 9556987
$ws.Range("J132").Value = 3176.8462
$ws.Range("K132").Value = 28670961
$ws.Range("L132").Value = 9530.5386
$ws.Range("M132").Value = -28668431
$ws.Range("N132").Value = -14590.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2400
$ws.Range("I96").Value = 2500
$ws.Range("K96").Value = 2500
$ws.Range("M96").Value = -1127

$ws.Range("H113").Value = 720.96
$ws.Range("I113").Value = 311
$ws.Range("J113").Value = 1242.7273
$ws.Range("K113").Value = 933
$ws.Range("L113").Value = 3728.1819
$ws.Range("M113").Value = 1237
$ws.Range("N113").Value = -8068.1819

$ws.Range("H126").Value = 989.5
$ws.Range("I126").Value = 743.125
$ws.Range("K126").Value = 2229.375
$ws.Range("M126").Value = 240.625

$ws.Range("H132").Value = 1186.1957
$ws.Range("I132").Value = 863.34283
$ws.Range("J132").Value = 2213.4546
$ws.Range("K132").Value = 2590.02849
$ws.Range("L132").Value = 6640.3638
$ws.Range("M132").Value = -60.02849000000015
$ws.Range("N132").Value = -11700.3638
